$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 45887
$ws.Range("B2").Value = 109.52
$ws.Range("C2").Value = 105.54
$ws.Range("D2").Value = 104.3
$ws.Range("E2").Value = 103.1
$ws.Range("F2").Value = 99.81999999999999
$ws.Range("G2").Value = 101.12
$ws.Range("H2").Value = 105.34
$ws.Range("I2").Value = 106.54
$ws.Range("J2").Value = 106.54
$ws.Range("K2").Value = 104.56
$ws.Range("L2").Value = 97.52
$ws.Range("M2").Value = 75.59
$ws.Range("N2").Value = 69
$ws.Range("O2").Value = 69
$ws.Range("P2").Value = 67.90000000000001
$ws.Range("Q2").Value = 60
$ws.Range("R2").Value = 69.26000000000001
$ws.Range("S2").Value = 84.53
$ws.Range("T2").Value = 104.07
$ws.Range("U2").Value = 115.06
$ws.Range("V2").Value = 120
$ws.Range("W2").Value = 120.76
$ws.Range("X2").Value = 114.32
$ws.Range("Y2").Value = 106.84
$ws.Range("Z2").Value = 96.68000000000001
$ws.Range("AB2").Value = 115.48
$ws.Range("AC2").Value = "20h-22h"
$ws.Range("AD2").Value = 120.38
$ws.Range("AE2").Value = "22h-24h"
$ws.Range("AF2").Value = 110.58
$ws.Range("AG2").Value = "11h-17h"
